$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -10
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = 2
